$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data to append at row 3, columns A:BD
$rowData = @(
    '0QZWAqio', '25/10/2024', '08:00', 'BULGARIA - PARVA LIGA', 'Arda', 'Septemvri Sofia',
    1.48, 4.1, 7, 2.05, 2.25, 7, 1.06, 10, 1.3, 3.4, 2.05, 1.8, 1.4, 2.75, 2.2, 1.62,
    6, 6.5, 9, 9.5, 13, 34, 9.5, 8, 21, 67, 201, 15, 34, 21, 81, 51, 51,
    3.25, 7.5, 21, 23, 51, 151, 2.75, 9.5, 67, 51, 8, 41, 41, 151, 151, 351, 51
)

for ($i = 0; $i -lt $rowData.Length; $i++) {
    $cell = $ws.Cells.Item(3, $i + 1)
    $cell.Value = $rowData[$i]
}
